$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12 updates
$ws.Range("D12").Value = [DateTime]"2022-06-13"
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = 8000
$ws.Range("N12").Value = '$/caja 36 atados'
$ws.Range("P12").Value = 222
$ws.Range("Q12").Value = 36
